# Slide 1 ("TextBox 3" / id=4) holds the "Presented By:" credit block.
# The author:
#   - reworded two of the lines ("Collage:" / "Branch:" prefixes)
#   - added a new trailing line with the AICTE student id
#   - grew the text box to fit the extra line
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 3")
$tr = $shape.TextFrame.TextRange

$cr = [char]13
$fullText = "Presented By:" + $cr + `
    "Raghav Joshi" + $cr + `
    "Collage: JIMS Engineering Management Technical Campus" + $cr + `
    "Branch: Computer Science Engineering" + $cr + `
    "AICTE Student ID: STU68427b55333691749187413"

# Re-assigning the whole TextRange in one shot (rather than touching each
# paragraph individually) lets PowerPoint re-flow/re-tag every run
# consistently, which is what keeps the existing bold/size/color formatting
# intact on the reworded + brand new lines alike.
$tr.Text = $fullText

# This is an auto-fit ("shrink/grow shape to fit text") text box, so adding
# a 5th line of 20pt bold text grows its height by one more line.
$shape.Height = 128.44221
